# Edit for sheet "week": break out stock.yaml completed
# - rows 587-612: convert column D (bsecode) from text to numeric
# - rows 613-638: newly appended rows (bsecode kept as text), refreshed as of 30/11/2024 18:35:49

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week")

# Convert bsecode (column D) for existing rows 587-612 from text to numeric
$ws.Cells.Item(587, 4).Value = 532466
$ws.Cells.Item(588, 4).Value = 532541
$ws.Cells.Item(589, 4).Value = 505200
$ws.Cells.Item(590, 4).Value = 540762
$ws.Cells.Item(591, 4).Value = 500495
$ws.Cells.Item(592, 4).Value = 532175
$ws.Cells.Item(593, 4).Value = 500271
$ws.Cells.Item(594, 4).Value = 500770
$ws.Cells.Item(595, 4).Value = 543220
$ws.Cells.Item(596, 4).Value = 539268
$ws.Cells.Item(597, 4).Value = 532508
$ws.Cells.Item(598, 4).Value = 543300
$ws.Cells.Item(599, 4).Value = 540777
$ws.Cells.Item(600, 4).Value = 500253
$ws.Cells.Item(601, 4).Value = 532814
$ws.Cells.Item(602, 4).Value = 532810
$ws.Cells.Item(603, 4).Value = 500049
$ws.Cells.Item(604, 4).Value = 500103
$ws.Cells.Item(605, 4).Value = 532210
$ws.Cells.Item(606, 4).Value = 531213
$ws.Cells.Item(607, 4).Value = 543257
$ws.Cells.Item(608, 4).Value = 500183
$ws.Cells.Item(609, 4).Value = 532477
$ws.Cells.Item(610, 4).Value = 532149
$ws.Cells.Item(611, 4).Value = 532461
$ws.Cells.Item(612, 4).Value = 533098

# Append new rows 613-638 (bsecode stays text for these newly-added rows)
$ws.Cells.Item(613, 1).Value = 1
$ws.Cells.Item(613, 2).Value = "OFSS"
$ws.Cells.Item(613, 3).Value = "Oracle Financial Services Software Limited"
$ws.Cells.Item(613, 4).NumberFormat = "@"
$ws.Cells.Item(613, 4).Value = "532466"
$ws.Cells.Item(613, 5).Value = 0.33
$ws.Cells.Item(613, 6).Value = 11696.45
$ws.Cells.Item(613, 7).Value = 72200
$ws.Cells.Item(613, 8).Value = "week"
$ws.Cells.Item(613, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(614, 1).Value = 2
$ws.Cells.Item(614, 2).Value = "COFORGE"
$ws.Cells.Item(614, 3).Value = "Coforge (Niit Tech)"
$ws.Cells.Item(614, 4).NumberFormat = "@"
$ws.Cells.Item(614, 4).Value = "532541"
$ws.Cells.Item(614, 5).Value = 0.29
$ws.Cells.Item(614, 6).Value = 8685.85
$ws.Cells.Item(614, 7).Value = 226945
$ws.Cells.Item(614, 8).Value = "week"
$ws.Cells.Item(614, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(615, 1).Value = 3
$ws.Cells.Item(615, 2).Value = "EICHERMOT"
$ws.Cells.Item(615, 3).Value = "Eicher Motors Limited"
$ws.Cells.Item(615, 4).NumberFormat = "@"
$ws.Cells.Item(615, 4).Value = "505200"
$ws.Cells.Item(615, 5).Value = 0.34
$ws.Cells.Item(615, 6).Value = 4831.85
$ws.Cells.Item(615, 7).Value = 407774
$ws.Cells.Item(615, 8).Value = "week"
$ws.Cells.Item(615, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(616, 1).Value = 4
$ws.Cells.Item(616, 2).Value = "TIINDIA"
$ws.Cells.Item(616, 3).Value = "Tube Investments of India Ltd"
$ws.Cells.Item(616, 4).NumberFormat = "@"
$ws.Cells.Item(616, 4).Value = "540762"
$ws.Cells.Item(616, 5).Value = 0.33
$ws.Cells.Item(616, 6).Value = 3590.55
$ws.Cells.Item(616, 7).Value = 353301
$ws.Cells.Item(616, 8).Value = "week"
$ws.Cells.Item(616, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(617, 1).Value = 5
$ws.Cells.Item(617, 2).Value = "ESCORTS"
$ws.Cells.Item(617, 3).Value = "Escorts Limited"
$ws.Cells.Item(617, 4).NumberFormat = "@"
$ws.Cells.Item(617, 4).Value = "500495"
$ws.Cells.Item(617, 5).Value = 1.15
$ws.Cells.Item(617, 6).Value = 3553.9
$ws.Cells.Item(617, 7).Value = 59622
$ws.Cells.Item(617, 8).Value = "week"
$ws.Cells.Item(617, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(618, 1).Value = 6
$ws.Cells.Item(618, 2).Value = "CYIENT"
$ws.Cells.Item(618, 3).Value = "Cyient Limited"
$ws.Cells.Item(618, 4).NumberFormat = "@"
$ws.Cells.Item(618, 4).Value = "532175"
$ws.Cells.Item(618, 5).Value = -0.23
$ws.Cells.Item(618, 6).Value = 1853.45
$ws.Cells.Item(618, 7).Value = 173616
$ws.Cells.Item(618, 8).Value = "week"
$ws.Cells.Item(618, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(619, 1).Value = 7
$ws.Cells.Item(619, 2).Value = "MFSL"
$ws.Cells.Item(619, 3).Value = "Max Financial Services Limited"
$ws.Cells.Item(619, 4).NumberFormat = "@"
$ws.Cells.Item(619, 4).Value = "500271"
$ws.Cells.Item(619, 5).Value = -0.5600000000000001
$ws.Cells.Item(619, 6).Value = 1133.95
$ws.Cells.Item(619, 7).Value = 966284
$ws.Cells.Item(619, 8).Value = "week"
$ws.Cells.Item(619, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(620, 1).Value = 8
$ws.Cells.Item(620, 2).Value = "TATACHEM"
$ws.Cells.Item(620, 3).Value = "Tata Chemicals Limited"
$ws.Cells.Item(620, 4).NumberFormat = "@"
$ws.Cells.Item(620, 4).Value = "500770"
$ws.Cells.Item(620, 5).Value = 0.72
$ws.Cells.Item(620, 6).Value = 1111.75
$ws.Cells.Item(620, 7).Value = 334105
$ws.Cells.Item(620, 8).Value = "week"
$ws.Cells.Item(620, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(621, 1).Value = 9
$ws.Cells.Item(621, 2).Value = "MAXHEALTH"
$ws.Cells.Item(621, 3).Value = "Max Healthcare Institute Ltd"
$ws.Cells.Item(621, 4).NumberFormat = "@"
$ws.Cells.Item(621, 4).Value = "543220"
$ws.Cells.Item(621, 5).Value = 0.3
$ws.Cells.Item(621, 6).Value = 979.75
$ws.Cells.Item(621, 7).Value = 2495600
$ws.Cells.Item(621, 8).Value = "week"
$ws.Cells.Item(621, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(622, 1).Value = 10
$ws.Cells.Item(622, 2).Value = "SYNGENE"
$ws.Cells.Item(622, 3).Value = "Syngene International Limited"
$ws.Cells.Item(622, 4).NumberFormat = "@"
$ws.Cells.Item(622, 4).Value = "539268"
$ws.Cells.Item(622, 5).Value = 2.61
$ws.Cells.Item(622, 6).Value = 940.8
$ws.Cells.Item(622, 7).Value = 1275641
$ws.Cells.Item(622, 8).Value = "week"
$ws.Cells.Item(622, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(623, 1).Value = 11
$ws.Cells.Item(623, 2).Value = "JSL"
$ws.Cells.Item(623, 3).Value = "Jindal Stainless Limited"
$ws.Cells.Item(623, 4).NumberFormat = "@"
$ws.Cells.Item(623, 4).Value = "532508"
$ws.Cells.Item(623, 5).Value = -1.84
$ws.Cells.Item(623, 6).Value = 683.2
$ws.Cells.Item(623, 7).Value = 690263
$ws.Cells.Item(623, 8).Value = "week"
$ws.Cells.Item(623, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(624, 1).Value = 12
$ws.Cells.Item(624, 2).Value = "SONACOMS"
$ws.Cells.Item(624, 3).Value = "Sona BLW Precision Forgings Ltd"
$ws.Cells.Item(624, 4).NumberFormat = "@"
$ws.Cells.Item(624, 4).Value = "543300"
$ws.Cells.Item(624, 5).Value = 1.04
$ws.Cells.Item(624, 6).Value = 671.6
$ws.Cells.Item(624, 7).Value = 1776675
$ws.Cells.Item(624, 8).Value = "week"
$ws.Cells.Item(624, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(625, 1).Value = 13
$ws.Cells.Item(625, 2).Value = "HDFCLIFE"
$ws.Cells.Item(625, 3).Value = "HDFC Life Insurance Company Ltd"
$ws.Cells.Item(625, 4).NumberFormat = "@"
$ws.Cells.Item(625, 4).Value = "540777"
$ws.Cells.Item(625, 5).Value = 0.01
$ws.Cells.Item(625, 6).Value = 657.75
$ws.Cells.Item(625, 7).Value = 5934406
$ws.Cells.Item(625, 8).Value = "week"
$ws.Cells.Item(625, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(626, 1).Value = 14
$ws.Cells.Item(626, 2).Value = "LICHSGFIN"
$ws.Cells.Item(626, 3).Value = "Lic Housing Finance Limited"
$ws.Cells.Item(626, 4).NumberFormat = "@"
$ws.Cells.Item(626, 4).Value = "500253"
$ws.Cells.Item(626, 5).Value = 0.65
$ws.Cells.Item(626, 6).Value = 638.8
$ws.Cells.Item(626, 7).Value = 831686
$ws.Cells.Item(626, 8).Value = "week"
$ws.Cells.Item(626, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(627, 1).Value = 15
$ws.Cells.Item(627, 2).Value = "INDIANB"
$ws.Cells.Item(627, 3).Value = "Indian Bank"
$ws.Cells.Item(627, 4).NumberFormat = "@"
$ws.Cells.Item(627, 4).Value = "532814"
$ws.Cells.Item(627, 5).Value = 0.22
$ws.Cells.Item(627, 6).Value = 574.3
$ws.Cells.Item(627, 7).Value = 2166791
$ws.Cells.Item(627, 8).Value = "week"
$ws.Cells.Item(627, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(628, 1).Value = 16
$ws.Cells.Item(628, 2).Value = "PFC"
$ws.Cells.Item(628, 3).Value = "Power Finance Corporation Limited"
$ws.Cells.Item(628, 4).NumberFormat = "@"
$ws.Cells.Item(628, 4).Value = "532810"
$ws.Cells.Item(628, 5).Value = 0.26
$ws.Cells.Item(628, 6).Value = 495.3
$ws.Cells.Item(628, 7).Value = 7808841
$ws.Cells.Item(628, 8).Value = "week"
$ws.Cells.Item(628, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(629, 1).Value = 17
$ws.Cells.Item(629, 2).Value = "BEL"
$ws.Cells.Item(629, 3).Value = "Bharat Electronics Limited"
$ws.Cells.Item(629, 4).NumberFormat = "@"
$ws.Cells.Item(629, 4).Value = "500049"
$ws.Cells.Item(629, 5).Value = 0.74
$ws.Cells.Item(629, 6).Value = 308
$ws.Cells.Item(629, 7).Value = 23241947
$ws.Cells.Item(629, 8).Value = "week"
$ws.Cells.Item(629, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(630, 1).Value = 18
$ws.Cells.Item(630, 2).Value = "BHEL"
$ws.Cells.Item(630, 3).Value = "Bharat Heavy Electricals Limited"
$ws.Cells.Item(630, 4).NumberFormat = "@"
$ws.Cells.Item(630, 4).Value = "500103"
$ws.Cells.Item(630, 5).Value = -0.55
$ws.Cells.Item(630, 6).Value = 251.09
$ws.Cells.Item(630, 7).Value = 8530849
$ws.Cells.Item(630, 8).Value = "week"
$ws.Cells.Item(630, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(631, 1).Value = 19
$ws.Cells.Item(631, 2).Value = "CUB"
$ws.Cells.Item(631, 3).Value = "City Union Bank Limited"
$ws.Cells.Item(631, 4).NumberFormat = "@"
$ws.Cells.Item(631, 4).Value = "532210"
$ws.Cells.Item(631, 5).Value = -0.08
$ws.Cells.Item(631, 6).Value = 179.53
$ws.Cells.Item(631, 7).Value = 1420556
$ws.Cells.Item(631, 8).Value = "week"
$ws.Cells.Item(631, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(632, 1).Value = 20
$ws.Cells.Item(632, 2).Value = "MANAPPURAM"
$ws.Cells.Item(632, 3).Value = "Manappuram Finance Limited"
$ws.Cells.Item(632, 4).NumberFormat = "@"
$ws.Cells.Item(632, 4).Value = "531213"
$ws.Cells.Item(632, 5).Value = -0.29
$ws.Cells.Item(632, 6).Value = 156.26
$ws.Cells.Item(632, 7).Value = 3356222
$ws.Cells.Item(632, 8).Value = "week"
$ws.Cells.Item(632, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(633, 1).Value = 21
$ws.Cells.Item(633, 2).Value = "IRFC"
$ws.Cells.Item(633, 3).Value = "Indian Railway Finance Corporation Ltd"
$ws.Cells.Item(633, 4).NumberFormat = "@"
$ws.Cells.Item(633, 4).Value = "543257"
$ws.Cells.Item(633, 5).Value = -2.65
$ws.Cells.Item(633, 6).Value = 149.34
$ws.Cells.Item(633, 7).Value = 18576066
$ws.Cells.Item(633, 8).Value = "week"
$ws.Cells.Item(633, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(634, 1).Value = 22
$ws.Cells.Item(634, 2).Value = "HFCL"
$ws.Cells.Item(634, 3).Value = "Himachal Futuristic Communications Limited"
$ws.Cells.Item(634, 4).NumberFormat = "@"
$ws.Cells.Item(634, 4).Value = "500183"
$ws.Cells.Item(634, 5).Value = -3.08
$ws.Cells.Item(634, 6).Value = 129.09
$ws.Cells.Item(634, 7).Value = 20074283
$ws.Cells.Item(634, 8).Value = "week"
$ws.Cells.Item(634, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(635, 1).Value = 23
$ws.Cells.Item(635, 2).Value = "UNIONBANK"
$ws.Cells.Item(635, 3).Value = "Union Bank Of India"
$ws.Cells.Item(635, 4).NumberFormat = "@"
$ws.Cells.Item(635, 4).Value = "532477"
$ws.Cells.Item(635, 5).Value = -0.06
$ws.Cells.Item(635, 6).Value = 121.62
$ws.Cells.Item(635, 7).Value = 23869952
$ws.Cells.Item(635, 8).Value = "week"
$ws.Cells.Item(635, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(636, 1).Value = 24
$ws.Cells.Item(636, 2).Value = "BANKINDIA"
$ws.Cells.Item(636, 3).Value = "Bank Of India"
$ws.Cells.Item(636, 4).NumberFormat = "@"
$ws.Cells.Item(636, 4).Value = "532149"
$ws.Cells.Item(636, 5).Value = -1.25
$ws.Cells.Item(636, 6).Value = 110.5
$ws.Cells.Item(636, 7).Value = 13029139
$ws.Cells.Item(636, 8).Value = "week"
$ws.Cells.Item(636, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(637, 1).Value = 25
$ws.Cells.Item(637, 2).Value = "PNB"
$ws.Cells.Item(637, 3).Value = "Punjab National Bank"
$ws.Cells.Item(637, 4).NumberFormat = "@"
$ws.Cells.Item(637, 4).Value = "532461"
$ws.Cells.Item(637, 5).Value = -1.31
$ws.Cells.Item(637, 6).Value = 104.9
$ws.Cells.Item(637, 7).Value = 30163680
$ws.Cells.Item(637, 8).Value = "week"
$ws.Cells.Item(637, 9).Value = "30/11/2024 18:35:49"

$ws.Cells.Item(638, 1).Value = 26
$ws.Cells.Item(638, 2).Value = "NHPC"
$ws.Cells.Item(638, 3).Value = "Nhpc Limited"
$ws.Cells.Item(638, 4).NumberFormat = "@"
$ws.Cells.Item(638, 4).Value = "533098"
$ws.Cells.Item(638, 5).Value = -2.48
$ws.Cells.Item(638, 6).Value = 81.44
$ws.Cells.Item(638, 7).Value = 25198619
$ws.Cells.Item(638, 8).Value = "week"
$ws.Cells.Item(638, 9).Value = "30/11/2024 18:35:49"
